# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" colours, used by the Notes Master
#   ppt/theme/theme2.xml  -> "Integral" colours, used by the Slide Master /
#                            the presentation itself
#
# The authored change swaps the two themes' colour schemes: the theme that
# backs the slide master/presentation (theme2.xml) ends up holding the
# default "Office Theme" palette, while the font scheme and format scheme
# (gradients, line styles, effects, …) - which were already identical
# between the two parts - stay untouched.
#
# Helper: turn an "RRGGBB" hex string into the BGR-packed long that
# PowerPoint's ColorFormat/ThemeColor .RGB property expects
# (PowerPoint stores/returns RGB() as R + G*256 + B*65536).
function ConvertTo-PptRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The presentation's theme colour scheme (12 slots, in the fixed OOXML
# order below) is reachable from the slide master.
$themeColors = $p.Slides.Item(1).Master.Theme.ThemeColorScheme

# Target palette: the stock "Office Theme" colours.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $themeColors.Item($i + 1).RGB = ConvertTo-PptRgb $officeThemeColors[$i]
}
